$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Add "Day 8" sheet (LeetCode "Users" table) after the last existing sheet.
# Cells are written header-row-first, then column by column (matching how
# the original workbook's pandas/openpyxl export ordered the shared strings)
# so the shared-string table comes out in the same order as the target file.
# ---------------------------------------------------------------------------
$day8 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$day8.Name = "Day 8"

$day8Headers = @("user_id", "name", "mail")
$day8Body = @(
    @(1, "Winston", "winston@leetcode.com"),
    @(2, "Jonathan", "jonathanisgreat"),
    @(3, "Annabelle", "bella-@leetcode.com"),
    @(4, "Sally", "sally.come@leetcode.com"),
    @(5, "Marwan", "quarz#2020@leetcode.com"),
    @(6, "David", "david69@gmail.com"),
    @(7, "Shapiro", ".shapo@leetcode.com")
)

for ($c = 0; $c -lt $day8Headers.Length; $c++) {
    $day8.Cells.Item(1, $c + 1).Value = $day8Headers[$c]
}
for ($c = 0; $c -lt $day8Headers.Length; $c++) {
    for ($r = 0; $r -lt $day8Body.Length; $r++) {
        $day8.Cells.Item($r + 2, $c + 1).Value = $day8Body[$r][$c]
    }
}

# ---------------------------------------------------------------------------
# Add "Day 9" sheet (LeetCode "Patients" table) after "Day 8"
# ---------------------------------------------------------------------------
$day9 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$day9.Name = "Day 9"

$day9Headers = @("patient_id", "patient_name", "conditions")
$day9Body = @(
    @(1, "Daniel", "YFEV COUGH"),
    @(2, "Alice", $null),
    @(3, "Bob", "DIAB100 MYOP"),
    @(4, "George", "ACNE DIAB100"),
    @(5, "Alain", "DIAB201")
)

for ($c = 0; $c -lt $day9Headers.Length; $c++) {
    $day9.Cells.Item(1, $c + 1).Value = $day9Headers[$c]
}
for ($c = 0; $c -lt $day9Headers.Length; $c++) {
    for ($r = 0; $r -lt $day9Body.Length; $r++) {
        $val = $day9Body[$r][$c]
        if ($null -ne $val) {
            $day9.Cells.Item($r + 2, $c + 1).Value = $val
        }
    }
}

# Match the best-fit column widths recorded for "Day 9" in the target file.
$day9.Columns.Item(1).ColumnWidth = 8.436197916666666
$day9.Columns.Item(2).ColumnWidth = 11.619791666666666
$day9.Columns.Item(3).ColumnWidth = 12.709635416666666

# ---------------------------------------------------------------------------
# View/selection tweaks on existing sheets
# ---------------------------------------------------------------------------

# "Day 6": rezoom and select the whole table, no longer the active tab
$day6 = $wb.Worksheets.Item("Day 6")
$day6.Activate()
$excel.ActiveWindow.Zoom = 189
$null = $day6.Range("A1:C6").Select()

# "Day 7": move the lingering selection
$day7 = $wb.Worksheets.Item("Day 7")
$day7.Activate()
$null = $day7.Range("C8").Select()

# "Day 9": leftover cursor position
$day9.Activate()
$null = $day9.Range("E10").Select()

# "Day 8" ends up the active/selected tab, matching the target workbook view
$day8.Activate()
$null = $day8.Range("I15").Select()
